# Corrects mis-scaled/incorrect financial figures in the IFRS company_list sheet
# (commit: 'error solve ifrs list'). Updates numeric cells D:AJ for rows 2-9;
# row 5 loses its V-column value (folded back into U5) and row 6 gains J6/O6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 10993
$ws.Range("E2").Value = 690
$ws.Range("F2").Value = 690
$ws.Range("G2").Value = 738
$ws.Range("H2").Value = 635
$ws.Range("I2").Value = 614
$ws.Range("J2").Value = 21
$ws.Range("K2").Value = 11150
$ws.Range("L2").Value = 4549
$ws.Range("M2").Value = 6602
$ws.Range("N2").Value = 5770
$ws.Range("O2").Value = 831
$ws.Range("P2").Value = 731
$ws.Range("Q2").Value = 438
$ws.Range("R2").Value = -135
$ws.Range("S2").Value = -561
$ws.Range("T2").Value = 337
$ws.Range("U2").Value = 101
$ws.Range("V2").Value = 255
$ws.Range("W2").Value = 6.28
$ws.Range("X2").Value = 5.78
$ws.Range("Y2").Value = 11.15
$ws.Range("Z2").Value = 5.47
$ws.Range("AA2").Value = 68.91
$ws.Range("AB2").Value = 657.2
$ws.Range("AC2").Value = 4198
$ws.Range("AD2").Value = 10.01
$ws.Range("AE2").Value = 40574
$ws.Range("AF2").Value = 1.04
$ws.Range("AG2").Value = 800
$ws.Range("AH2").Value = 1.9
$ws.Range("AI2").Value = 18.54
$ws.Range("AJ2").Value = 14403386

# Row 3
$ws.Range("D3").Value = 12105
$ws.Range("E3").Value = 1231
$ws.Range("F3").Value = 1231
$ws.Range("G3").Value = 1327
$ws.Range("H3").Value = 974
$ws.Range("I3").Value = 856
$ws.Range("J3").Value = 118
$ws.Range("K3").Value = 11728
$ws.Range("L3").Value = 4039
$ws.Range("M3").Value = 7689
$ws.Range("N3").Value = 6599
$ws.Range("O3").Value = 1089
$ws.Range("P3").Value = 731
$ws.Range("Q3").Value = 705
$ws.Range("R3").Value = -124
$ws.Range("S3").Value = -100
$ws.Range("T3").Value = 264
$ws.Range("U3").Value = 441
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 10.17
$ws.Range("X3").Value = 8.039999999999999
$ws.Range("Y3").Value = 13.84
$ws.Range("Z3").Value = 8.51
$ws.Range("AA3").Value = 52.53
$ws.Range("AB3").Value = 767.71
$ws.Range("AC3").Value = 5852
$ws.Range("AD3").Value = 14.66
$ws.Range("AE3").Value = 45341
$ws.Range("AF3").Value = 1.89
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1.17
$ws.Range("AI3").Value = 17.01
$ws.Range("AJ3").Value = 14403386

# Row 4
$ws.Range("D4").Value = 11538
$ws.Range("E4").Value = 716
$ws.Range("F4").Value = 716
$ws.Range("G4").Value = 603
$ws.Range("H4").Value = 343
$ws.Range("I4").Value = 336
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 12482
$ws.Range("L4").Value = 4594
$ws.Range("M4").Value = 7888
$ws.Range("N4").Value = 6806
$ws.Range("O4").Value = 1082
$ws.Range("P4").Value = 731
$ws.Range("Q4").Value = 1061
$ws.Range("R4").Value = -354
$ws.Range("S4").Value = -171
$ws.Range("T4").Value = 395
$ws.Range("U4").Value = 666
$ws.Range("V4").Value = 3
$ws.Range("W4").Value = 6.21
$ws.Range("X4").Value = 2.98
$ws.Range("Y4").Value = 5.01
$ws.Range("Z4").Value = 2.84
$ws.Range("AA4").Value = 58.24
$ws.Range("AB4").Value = 799.85
$ws.Range("AC4").Value = 2295
$ws.Range("AD4").Value = 21.28
$ws.Range("AE4").Value = 46542
$ws.Range("AF4").Value = 1.05
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 2.05
$ws.Range("AI4").Value = 43.57
$ws.Range("AJ4").Value = 14623136

# Row 5
$ws.Range("D5").Value = 11375
$ws.Range("E5").Value = 967
$ws.Range("F5").Value = 967
$ws.Range("G5").Value = 720
$ws.Range("H5").Value = 591
$ws.Range("I5").Value = 569
$ws.Range("J5").Value = 22
$ws.Range("K5").Value = 12263
$ws.Range("L5").Value = 3975
$ws.Range("M5").Value = 8288
$ws.Range("N5").Value = 7205
$ws.Range("O5").Value = 1083
$ws.Range("P5").Value = 731
$ws.Range("Q5").Value = 703
$ws.Range("R5").Value = -206
$ws.Range("S5").Value = -157
$ws.Range("T5").Value = 187
$ws.Range("U5").Value = 516
$ws.Range("V5").ClearContents()
$ws.Range("W5").Value = 8.51
$ws.Range("X5").Value = 5.19
$ws.Range("Y5").Value = 8.119999999999999
$ws.Range("Z5").Value = 4.77
$ws.Range("AA5").Value = 47.97
$ws.Range("AB5").Value = 857.16
$ws.Range("AC5").Value = 3888
$ws.Range("AD5").Value = 12.05
$ws.Range("AE5").Value = 49272
$ws.Range("AF5").Value = 0.95
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 2.13
$ws.Range("AI5").Value = 25.72
$ws.Range("AJ5").Value = 14623136

# Row 6
$ws.Range("D6").Value = 10584
$ws.Range("E6").Value = 563
$ws.Range("F6").Value = 563
$ws.Range("G6").Value = 765
$ws.Range("H6").Value = 572
$ws.Range("I6").Value = 563
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 11965
$ws.Range("L6").Value = 3717
$ws.Range("M6").Value = 8248
$ws.Range("N6").Value = 7276
$ws.Range("O6").Value = 972
$ws.Range("P6").Value = 731
$ws.Range("Q6").Value = 904
$ws.Range("R6").Value = -609
$ws.Range("S6").Value = -539
$ws.Range("T6").Value = 271
$ws.Range("U6").Value = 632
$ws.Range("V6").Value = 65
$ws.Range("W6").Value = 5.32
$ws.Range("X6").Value = 5.41
$ws.Range("Y6").Value = 7.77
$ws.Range("Z6").Value = 4.72
$ws.Range("AA6").Value = 45.07
$ws.Range("AB6").Value = 921.3099999999999
$ws.Range("AC6").Value = 3848
$ws.Range("AD6").Value = 7.03
$ws.Range("AE6").Value = 54639
$ws.Range("AF6").Value = 0.5
$ws.Range("AG6").Value = 1200
$ws.Range("AH6").Value = 4.44
$ws.Range("AI6").Value = 28.4
$ws.Range("AJ6").Value = 14623136

# Row 7
$ws.Range("D7").Value = 9939
$ws.Range("E7").Value = 888
$ws.Range("G7").Value = 1059
$ws.Range("H7").Value = 75
$ws.Range("I7").Value = 65
$ws.Range("K7").Value = 11429
$ws.Range("L7").Value = 3507
$ws.Range("M7").Value = 7922
$ws.Range("N7").Value = 7163
$ws.Range("P7").Value = 730
$ws.Range("Q7").Value = 461
$ws.Range("R7").Value = -254
$ws.Range("S7").Value = -178
$ws.Range("T7").Value = 322
$ws.Range("U7").Value = 444
$ws.Range("W7").Value = 8.93
$ws.Range("X7").Value = 0.75
$ws.Range("Y7").Value = 0.9
$ws.Range("Z7").Value = 0.64
$ws.Range("AA7").Value = 44.27
$ws.Range("AC7").Value = 444
$ws.Range("AD7").Value = 91.11
$ws.Range("AE7").Value = 51825
$ws.Range("AF7").Value = 0.78
$ws.Range("AG7").Value = 1233
$ws.Range("AH7").Value = 3.05
$ws.Range("AI7").Value = 277.46

# Row 8
$ws.Range("D8").Value = 10890
$ws.Range("E8").Value = 880
$ws.Range("G8").Value = 965
$ws.Range("H8").Value = 745
$ws.Range("I8").Value = 740
$ws.Range("K8").Value = 12560
$ws.Range("L8").Value = 3840
$ws.Range("M8").Value = 8720
$ws.Range("N8").Value = 7735
$ws.Range("P8").Value = 730
$ws.Range("Q8").Value = 905
$ws.Range("R8").Value = -410
$ws.Range("S8").Value = -195
$ws.Range("T8").Value = 365
$ws.Range("U8").Value = 910
$ws.Range("W8").Value = 8.08
$ws.Range("X8").Value = 6.84
$ws.Range("Y8").Value = 9.92
$ws.Range("Z8").Value = 6.12
$ws.Range("AA8").Value = 44.04
$ws.Range("AC8").Value = 5060
$ws.Range("AD8").Value = 7.7
$ws.Range("AE8").Value = 55964
$ws.Range("AF8").Value = 0.7
$ws.Range("AG8").Value = 1450
$ws.Range("AH8").Value = 3.72
$ws.Range("AI8").Value = 28.65

# Row 9
$ws.Range("D9").Value = 12335
$ws.Range("E9").Value = 970
$ws.Range("G9").Value = 1115
$ws.Range("H9").Value = 860
$ws.Range("I9").Value = 850
$ws.Range("K9").Value = 13725
$ws.Range("L9").Value = 4345
$ws.Range("M9").Value = 9375
$ws.Range("N9").Value = 8390
$ws.Range("P9").Value = 730
$ws.Range("Q9").Value = 1145
$ws.Range("R9").Value = -460
$ws.Range("S9").Value = -210
$ws.Range("T9").Value = 370
$ws.Range("U9").Value = 1020
$ws.Range("W9").Value = 7.86
$ws.Range("X9").Value = 6.97
$ws.Range("Y9").Value = 10.54
$ws.Range("Z9").Value = 6.54
$ws.Range("AA9").Value = 46.35
$ws.Range("AC9").Value = 5813
$ws.Range("AD9").Value = 6.7
$ws.Range("AE9").Value = 60703
$ws.Range("AF9").Value = 0.64
$ws.Range("AG9").Value = 1500
$ws.Range("AH9").Value = 3.85
$ws.Range("AI9").Value = 25.81

Write-Host "edits applied"
